$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.466.18"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "1.667.55"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.29"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.00"
$ws.Range("E9").Value = "  +5.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.407"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9987"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08585"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.47"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.327"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001342"
$ws.Range("E15").Value = "  +3.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.885"
$ws.Range("E16").Value = "  +4.30%  "
$ws.Range("D17").Value = "1.657.22"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.58"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06967"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.56"
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.000"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9980"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.72"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").Value = "24.465.96"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.435"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.016"
$ws.Range("E26").Value = "  +8.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.53"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.90"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "142.75"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.384"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.084"
$ws.Range("E31").Value = "  -7.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.507"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").Value = "1.839.88"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.065"
$ws.Range("E34").Value = "  +6.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08259"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.28"
$ws.Range("E36").Value = "  +11.34%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02987"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.790"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2756"
$ws.Range("E39").Value = "  +0.98%  "
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7739"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.85"
$ws.Range("E42").Value = "  +4.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.446"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.47"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7110"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.530"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.143"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9985"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08456"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.24"
$ws.Range("E51").Value = "  +12.03%  "
